# Insert a new data row before the current row 27 (shifts rows 27-55 down to 28-56)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new record.
# Columns A,B,C,E,F,G,H,I,O,Q,R mirror the record that used to sit at row 27
# (now pushed down to row 28); only D, J, K, L, M, N, P differ.
$ws.Range("A27").Value = 3
$ws.Range("B27").Value = "Femacal de La Calera"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value = 44629
$ws.Range("E27").Value = 5
$ws.Range("F27").Value = 100112022
$ws.Range("G27").Value = "Arveja Verde"
$ws.Range("H27").Value = "Perfection"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 45
$ws.Range("K27").Value = 24000
$ws.Range("L27").Value = 25000
$ws.Range("M27").Value = 24444
$ws.Range("N27").Value = "$/saco 25 kilos"
$ws.Range("O27").Value = "Región Metropolitana"
$ws.Range("P27").Value = 978
$ws.Range("Q27").Value = 25
$ws.Range("R27").Value = "Hortaliza"

# Match the date display/number format used by the rest of column D.
$ws.Range("D27").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "Row inserted at 27"
